# Update the "Förändrad" (Changed) date column (C) for rows 2-13
# from serial date 45243 (2023-11-13) to 45244 (2023-11-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Range("C$row")
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
